$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be parsed as numbers
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "69.424.73"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "3.686.66"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "685.23"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "160.08"
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("D7").Value = "3.685.40"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -4.74%  "
$ws.Range("D10").Value = "0.146"
$ws.Range("E10").Value = "  -7.58%  "
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").Value = "  -4.79%  "
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  -7.92%  "
$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").Value = "  -4.99%  "
$ws.Range("D14").Value = "4.305.74"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "32.52"
$ws.Range("E15").Value = "  -8.57%  "
$ws.Range("D16").Value = "3.683.30"
$ws.Range("E16").Value = "  -4.46%  "
$ws.Range("D17").Value = "69.447.53"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "15.86"
$ws.Range("E19").Value = "  -8.29%  "
$ws.Range("D20").Value = "6.46"
$ws.Range("E20").Value = "  -8.81%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "10.18"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "469.64"
$ws.Range("E22").Value = "  -8.55%  "
$ws.Range("D23").Value = "0.648"
$ws.Range("E23").Value = "  -8.31%  "
$ws.Range("D24").Value = "79.67"
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("D25").Value = "3.831.66"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -10.56%  "
$ws.Range("D28").Value = "10.96"
$ws.Range("E28").Value = "  -11.72%  "
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  -8.77%  "
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").Value = "  -7.25%  "
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  -9.60%  "
$ws.Range("D32").Value = "2.01"
$ws.Range("E32").Value = "  -8.93%  "
$ws.Range("D33").Value = "6.60"
$ws.Range("E33").Value = "  -8.55%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "26.76"
$ws.Range("E35").Value = "  -6.97%  "
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.659.88"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "8.17"
$ws.Range("E38").Value = "  -10.70%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "6.13"
$ws.Range("E39").Value = "  -5.60%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0902"
$ws.Range("E42").Value = "  -8.36%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "167.00"
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.943"
$ws.Range("E45").Value = "  -5.79%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "47.62"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("B47").Value = "SuiNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D47").Value = "1.14"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  -11.74%  "
$ws.Range("D49").Value = "0.000278"
$ws.Range("E49").Value = "  -5.88%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.30"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "28.37"
$ws.Range("E51").Value = "  -1.84%  "
